# Set max ker degrowth
# Edits the RSD sheet of the workbook:
#  1. Insert a new row at row 17 (this pushes the existing "Solar" growth-rate
#     row down to row 18, and also cascades the shift through the data table
#     below, duplicating the last "Solar" data row into a new row 33).
#  2. Populate the newly inserted row 17 with a "MaxDegrowth" unit constraint
#     for Kerosene (mirrors the MaxGrowth row above it, but subtracts instead
#     of adds the rate, and points at the Kerosene data row).
#  3. Halve the Gas max-growth-rate input (now on row 29) from 0.1 to 0.05.
#  4. Make the RSD sheet the active tab/selection, matching the saved view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RSD")

# 1. Insert a new row before row 17 (shifts row 17 onward down by one).
$ws.Rows.Item(17).Insert()

# 2. Fill in the new row 17 - Kerosene "MaxDegrowth" unit constraint.
$ws.Range("B17").Formula = '=TEXTJOIN("_",TRUE,"UC",A32,"MaxDegrowth",B32)'
$ws.Range("C17").Value = "ACT, GROWTH"
$ws.Range("F17").Value = "RSDKER"
$ws.Range("G17").Value = "FT*"
$ws.Range("H17").Value = 2021
$ws.Range("I17").Value = "LO"
$ws.Range("J17").Formula = "=1-C33"
$ws.Range("K17").Value = 1
$ws.Range("L17").Formula = "=-D33"
$ws.Range("M17").Value = 5
$ws.Range("N17").Formula = '=TEXTJOIN(" ",TRUE,A32, "maximum degrowth rate of",B32)'

# 3. Update the Gas max growth rate input (data table row shifted to 29).
$ws.Range("C29").Value = 0.05

# 4. Select the RSD sheet / cell as the active view.
$ws.Activate()
$ws.Range("I26").Select()
